$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the lone row that held only the professor's name in columns B/C).
# This shifts every following row up by one, matching the workbook's new
# dimension of A1:C24 (one row fewer than before).
$ws.Rows(13).Delete()

# After the shift, a handful of cells end up holding different text than a
# plain shift would produce, so set their final content explicitly.

# Row 10 ("Objetivos:" answer) now holds the professor's name.
$ws.Range("B10:C10").Value = '6270264 - Juan Fernando Zapata Zapata'

# Row 13 ("Programa resumido:" answer) now holds "Semestral".
$ws.Range("B13:C13").Value = 'Semestral'

# Row 15 ("Programa:" answer) now holds the activation date. Build it as a
# formula result in a scratch cell and paste only the value across, so Excel
# keeps it as literal text instead of auto-converting "01/01/2018" into a
# date value (and without mutating any cell styles).
$ws.Range("A100").Formula = '="01/01/2018"'
$ws.Range("A100").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("A100").Clear()

# Row 18 ("Método:" answer) now holds the professor's name again.
$ws.Range("B18:C18").Value = '6270264 - Juan Fernando Zapata Zapata'

# Row 19 ("Critério:" answer) now holds the evaluation method text.
$ws.Range("B19:C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# Row 20 ("Norma de recuperação:" answer) now holds the NF>=5,0 criterion text.
$ws.Range("B20:C20").Value = 'NF≥ 5,0.'

# Row 21 ("Bibliografia:" answer) now holds the recovery norm text.
$ws.Range("B21:C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
